$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new rows before row 111; existing rows 111-140 shift down to 113-142
$ws.Rows.Item(111).Insert()
$ws.Rows.Item(111).Insert()

# New row 111: Murcott / Primera, Region de O'Higgins
$ws.Cells.Item(111, 1).Value = 11
$ws.Cells.Item(111, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(111, 3).Value = "Bíobío"
$ws.Cells.Item(111, 4).Value = 44798
$ws.Cells.Item(111, 5).Value = 8
$ws.Cells.Item(111, 6).Value = "Fruta"
$ws.Cells.Item(111, 7).Value = 100102
$ws.Cells.Item(111, 8).Value = "Cítricos"
$ws.Cells.Item(111, 9).Value = 100102004
$ws.Cells.Item(111, 10).Value = "Mandarina"
$ws.Cells.Item(111, 11).Value = "Murcott"
$ws.Cells.Item(111, 12).Value = "Primera"
$ws.Cells.Item(111, 13).Value = 290
$ws.Cells.Item(111, 14).Value = 8000
$ws.Cells.Item(111, 15).Value = 8500
$ws.Cells.Item(111, 16).Value = 8241
$ws.Cells.Item(111, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(111, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(111, 19).Value = 458
$ws.Cells.Item(111, 20).Value = 18

# New row 112: Murcott / Segunda, Region de O'Higgins
$ws.Cells.Item(112, 1).Value = 11
$ws.Cells.Item(112, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(112, 3).Value = "Bíobío"
$ws.Cells.Item(112, 4).Value = 44798
$ws.Cells.Item(112, 5).Value = 8
$ws.Cells.Item(112, 6).Value = "Fruta"
$ws.Cells.Item(112, 7).Value = 100102
$ws.Cells.Item(112, 8).Value = "Cítricos"
$ws.Cells.Item(112, 9).Value = 100102004
$ws.Cells.Item(112, 10).Value = "Mandarina"
$ws.Cells.Item(112, 11).Value = "Murcott"
$ws.Cells.Item(112, 12).Value = "Segunda"
$ws.Cells.Item(112, 13).Value = 220
$ws.Cells.Item(112, 14).Value = 6500
$ws.Cells.Item(112, 15).Value = 6500
$ws.Cells.Item(112, 16).Value = 6273
$ws.Cells.Item(112, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(112, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(112, 19).Value = 348
$ws.Cells.Item(112, 20).Value = 18
